$d = $word.ActiveDocument
$d.Content.Find.Execute("Distance from the school to the trust headquarters", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Distance from the converting school to the trust, or other schools in the trust", 2)
